# Append two new heading/sub-heading paragraphs at the very end of the
# document body (after the last, already-empty paragraph that follows the
# "Revision Summary" table, and before the section properties).
#
# Strategy: Word's Range.InsertXML() merges the LAST <w:p> of the supplied
# OOXML fragment into the paragraph the (collapsed) range currently sits
# in -- preserving that host paragraph's own pPr/identity -- while any
# preceding <w:p> elements in the fragment are inserted verbatim as brand
# new sibling paragraphs (full control over their own pPr/rPr). To keep
# the existing trailing empty paragraph completely untouched, we first
# create a fresh scratch paragraph after it (InsertParagraphAfter) to act
# as the merge "host", insert our two fully-formed paragraphs plus a
# trailing empty dummy paragraph (which absorbs the merge into the host),
# and finally delete that now-redundant host paragraph.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Create a scratch paragraph right after the existing trailing empty one;
# this will be the "host" that InsertXML merges its final fragment into.
$lastPara.Range.InsertParagraphAfter()

$hostIndex = $d.Paragraphs.Count
$hostPara = $d.Paragraphs.Item($hostIndex)
$insertionPoint = $d.Range($hostPara.Range.Start, $hostPara.Range.Start)

$para1 = "<w:p $wNs>" +
  "<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/><w:highlight w:val=`"black`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"FFFFFF`" w:themeColor=`"background1`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/><w:highlight w:val=`"black`"/></w:rPr><w:t>1 Feature Descript</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"FFFFFF`" w:themeColor=`"background1`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/><w:highlight w:val=`"black`"/></w:rPr><w:t xml:space=`"preserve`">ion                                                                                                           </w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/><w:highlight w:val=`"black`"/></w:rPr><w:t>I</w:t></w:r>" +
  "</w:p>"

$para2 = "<w:p $wNs>" +
  "<w:pPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:color w:val=`"FFFFFF`" w:themeColor=`"background1`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>1</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>1.1</w:t></w:r>" +
  "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`"> Behavioral Specification</w:t></w:r>" +
  "</w:p>"

# Trailing empty dummy paragraph (with an explicit empty run) that absorbs
# the merge into the scratch host paragraph, leaving it a no-op.
$dummy = "<w:p $wNs><w:r></w:r></w:p>"

$insertionPoint.InsertXML($para1 + $para2 + $dummy)

# The scratch host paragraph is now the very last paragraph in the body;
# delete it so only our two real paragraphs remain before the sectPr.
$finalCount = $d.Paragraphs.Count
$hostNow = $d.Paragraphs.Item($finalCount)
$hostRange = $d.Range($hostNow.Range.Start, $hostNow.Range.End)
$hostRange.Delete()

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
